# Weekly update: a new week's price record is added for
# "Terminal La Palmera de La Serena - Zanahoria", inserted at the top of the
# date-ordered block (row 185). This shifts the existing rows 185-225 down to
# 186-226 (carrying their data/formatting with them), and the new row 185
# receives the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 185; Excel shifts rows 185-225 down
# to 186-226 automatically, preserving cell styles (e.g. the date format on
# column D).
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with the new week's record. The
# descriptive columns (A,B,C,E,F,G,H,I,N,O,Q,R) are constant across this
# whole block, so reuse the same values as the rest of the series.
$ws.Cells.Item(185, 1).Value = 8
$ws.Cells.Item(185, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(185, 3).Value = "Coquimbo"
$ws.Cells.Item(185, 4).Value = 44511
$ws.Cells.Item(185, 5).Value = 4
$ws.Cells.Item(185, 6).Value = 100114013
$ws.Cells.Item(185, 7).Value = "Zanahoria"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Primera"
$ws.Cells.Item(185, 10).Value = 700
$ws.Cells.Item(185, 11).Value = 6300
$ws.Cells.Item(185, 12).Value = 6500
$ws.Cells.Item(185, 13).Value = 6400
$ws.Cells.Item(185, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(185, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(185, 16).Value = 320
$ws.Cells.Item(185, 17).Value = 20
$ws.Cells.Item(185, 18).Value = "Hortaliza"
